$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Add new row 12 data (mirrors the new <row r="12"> in Logs sheet)
$ws.Range("A12").Value = "Geen onderwerp"
$ws.Range("B12").Value = "onbekend"
$ws.Range("D12").Value = "Overig"
$ws.Range("F12").Value = "2025-08-18 21:35:03"
$ws.Range("G12").Value = "Nee"
$ws.Range("H12").Value = "Ja"
$ws.Range("I12").Value = "Nee"
$ws.Range("J12").Value = "Nee"

# Extend conditional formatting ranges from row 11 to row 12
$colsToExtend = @("D", "G", "H", "I", "J")
foreach ($col in $colsToExtend) {
    $oldRange = $ws.Range($col + "2:" + $col + "11")
    $newRange = $ws.Range($col + "2:" + $col + "12")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# Update Dashboard sheet count for "Overig" from 3 to 4
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B3").Value = 4
